$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the existing "Denis" value in column C, row 8
$ws.Range("C8").Value = "Denis"

# Continue filling column C downward with new names
$ws.Range("C10").Value = "Hao"
$ws.Range("C11").Value = "Arbinnav"
$ws.Range("C12").Value = "Denvendra"
$ws.Range("C9").Value = "Gabriel"

# Shift the old A10:A11 values ("Yvoone", "See Fu") right into B10:B11
# to make room for a new name at the top of column A
$ws.Range("A10:A11").Cut($ws.Range("B10:B11"))
$ws.Range("A10").Value = "Faye"

# Update selection to match the resulting workbook state
$ws.Range("A11").Select()
